# B1--and-B2-PowerPoint.pptx
#
# 1) The table on slide 5 gets re-styled from the deck's custom "Table_0"
#    style to the built-in "Medium Style 2 - Accent 1" table style.
# 2) The presentation's theme colours are switched from the "Integral"
#    (Red Violet) palette to the stock "Office Theme" (Office) palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
foreach ($shp in $slide.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{83C4B8F4-E067-4A26-95A2-2FE4064CCCCC}")
    }
}

# --- 2) Swap the design's colour scheme to the Office palette -------------
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
